$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "mean_10Be_data_holocene"

# Row 2 - Siple Dome
$ws.Range("E2").Value = 14.2
$ws.Range("F2").Value = 107.16641298833079

# Row 3 - WD -> Wais Divide
$ws.Range("A3").Value = "Wais Divide"
$ws.Range("E3").Value = 20.1
$ws.Range("F3").Value = 91.143455098934552

# Row 4 - South Pole -> PS1
$ws.Range("A4").Value = "PS1"
$ws.Range("E4").Value = 8.72
$ws.Range("F4").Value = 107.83866057838662

# Row 5 - EDML
$ws.Range("E5").Value = 6.28
$ws.Range("F5").Value = 79.455859969558603

# Row 6 - EDC
$ws.Range("B6").Value = -75.06
$ws.Range("C6").Value = 123.21
$ws.Range("E6").Value = 2.5
$ws.Range("F6").Value = 53.67

# Row 7 - LDC
$ws.Range("D7").Value = 0.59
$ws.Range("E7").Value = 2.31
$ws.Range("F7").Value = 43.217275494672755

# Row 8 - ODC
$ws.Range("B8").Value = -74.39
$ws.Range("C8").Value = 124.1
$ws.Range("E8").Value = 2.99
$ws.Range("F8").Value = 48.354261796042621

# Row 9 - Vostok
$ws.Range("E9").Value = 1.93
$ws.Range("F9").Value = 53.917110603754438

# Row 10 - Dome F
$ws.Range("E10").Value = 3.21
$ws.Range("F10").Value = 88.555936073059357

# Row 11 - Dome A
$ws.Range("E11").Value = 2.29
$ws.Range("F11").Value = 73.486808726534761

# Row 12 - SPICE
$ws.Range("E12").Value = 8.16
$ws.Range("F12").Value = 90.304414003044144

# Restore the active selection cell
$ws.Range("F7").Select()
